$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2233
$ws.Range("J40").Value = 2549.5
$ws.Range("L40").Value = 2549.5
$ws.Range("N40").Value = -2899.5
# Row 92
$ws.Range("H92").Value = 2203.0454
$ws.Range("I92").Value = 2342.4375
$ws.Range("K92").Value = 2342.4375
$ws.Range("M92").Value = -1094.4375
# Row 137
$ws.Range("H137").Value = 40819704
$ws.Range("I137").Value = 28574330
$ws.Range("K137").Value = 85722990
$ws.Range("M137").Value = -85720440
# Row 138
$ws.Range("H138").Value = 6543126.5
$ws.Range("I138").Value = 3265.3845
$ws.Range("J138").Value = 8780448
$ws.Range("K138").Value = 9796.1535
$ws.Range("L138").Value = 26341344
$ws.Range("M138").Value = -4656.1535
$ws.Range("N138").Value = -26351624

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 125007064
$ws.Range("I61").Value = 333334660
$ws.Range("K61").Value = 333334660
$ws.Range("M61").Value = -333334448
# Row 102
$ws.Range("H102").Value = 3238
$ws.Range("I102").Value = 3112.1428
$ws.Range("K102").Value = 3112.1428
$ws.Range("M102").Value = -1490.1428
# Row 121
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = ""
# Row 132
$ws.Range("H132").Value = 41669920
$ws.Range("I132").Value = 3411.35
$ws.Range("J132").Value = 250002480
$ws.Range("K132").Value = 10234.05
$ws.Range("L132").Value = 750007440
$ws.Range("M132").Value = -7704.049999999999
$ws.Range("N132").Value = -750012500
# Row 136
$ws.Range("H136").Value = 125007064
$ws.Range("I136").Value = 333334660
$ws.Range("K136").Value = 1000003980
$ws.Range("M136").Value = -1000001430
# Row 139
$ws.Range("H139").Value = 61054.832
$ws.Range("J139").Value = 61054.832
$ws.Range("L139").Value = 61054.832
$ws.Range("N139").Value = -71334.83199999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 5769.75
$ws.Range("I99").Value = 5146.75
$ws.Range("K99").Value = 5146.75
$ws.Range("M99").Value = -3648.75
# Row 137
$ws.Range("H137").Value = 193311.8
$ws.Range("J137").Value = 193311.8
$ws.Range("L137").Value = 193311.8
$ws.Range("N137").Value = -203511.8
# Row 141
$ws.Range("H141").Value = 88730
$ws.Range("J141").Value = 88730
$ws.Range("L141").Value = 88730
$ws.Range("N141").Value = -99090

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 325.45456
$ws.Range("I7").Value = 103.8
$ws.Range("K7").Value = 103.8
$ws.Range("M7").Value = 9.200000000000003
# Row 16
$ws.Range("H16").Value = 2828.375
$ws.Range("I16").Value = 3217.8
$ws.Range("K16").Value = 3217.8
$ws.Range("M16").Value = -2930.8
# Row 31
$ws.Range("H31").Value = 33339584
$ws.Range("I31").Value = 5593.077
$ws.Range("J31").Value = 58830284
$ws.Range("K31").Value = 5593.077
$ws.Range("L31").Value = 58830284
$ws.Range("M31").Value = -5298.077
$ws.Range("N31").Value = -58830874
# Row 34
$ws.Range("H34").Value = 33339584
$ws.Range("I34").Value = 5593.077
$ws.Range("J34").Value = 58830284
$ws.Range("K34").Value = 5593.077
$ws.Range("L34").Value = 58830284
$ws.Range("M34").Value = -5391.077
$ws.Range("N34").Value = -58830688
# Row 52
$ws.Range("H52").Value = 149985.33
$ws.Range("J52").Value = 149985.33
$ws.Range("L52").Value = 149985.33
$ws.Range("N52").Value = -150573.33
# Row 113
$ws.Range("H113").Value = 2828.375
$ws.Range("I113").Value = 3217.8
$ws.Range("K113").Value = 3217.8
$ws.Range("M113").Value = -1047.8
# Row 131
$ws.Range("H131").Value = 89999
$ws.Range("J131").Value = 89999
$ws.Range("L131").Value = 89999
$ws.Range("N131").Value = -100079
# Row 141
$ws.Range("H141").Value = 288570.9
$ws.Range("J141").Value = 297443.06
$ws.Range("L141").Value = 297443.06
$ws.Range("N141").Value = -307803.06

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 1665671.5
$ws.Range("I2").Value = 1668
$ws.Range("J2").Value = 2422036.8
$ws.Range("K2").Value = 10008
$ws.Range("L2").Value = 14532220.8
$ws.Range("M2").Value = -9895
$ws.Range("N2").Value = -14532446.8
# Row 133
$ws.Range("H133").Value = 12690.1
$ws.Range("I133").Value = 7847
$ws.Range("K133").Value = 23541
$ws.Range("M133").Value = -18481
# Row 134
$ws.Range("H134").Value = 3588.353
$ws.Range("I134").Value = 1466.8667
$ws.Range("K134").Value = 4400.6001
$ws.Range("M134").Value = 669.3999000000003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 108.625
$ws.Range("I2").Value = 80
$ws.Range("K2").Value = 80
$ws.Range("M2").Value = 33
# Row 132
$ws.Range("H132").Value = 4942.0557
$ws.Range("I132").Value = 3896.1538
$ws.Range("K132").Value = 11688.4614
$ws.Range("M132").Value = -9158.4614
# Row 135
$ws.Range("H135").Value = 86838.94500000001
$ws.Range("J135").Value = 86838.94500000001
$ws.Range("L135").Value = 86838.94500000001
$ws.Range("N135").Value = -96978.94500000001
# Row 136
$ws.Range("H136").Value = 80000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 80000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 240000
$ws.Range("M136").Value = ""
$ws.Range("N136").Value = -245100

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4574.694
$ws.Range("J7").Value = 4840.143
$ws.Range("L7").Value = 4840.143
$ws.Range("N7").Value = -5064.143
# Row 22
$ws.Range("H22").Value = 3082.8333
$ws.Range("I22").Value = 1820
$ws.Range("J22").Value = 3984.8572
$ws.Range("K22").Value = 1820
$ws.Range("L22").Value = 3984.8572
$ws.Range("M22").Value = -1525
$ws.Range("N22").Value = -4574.8572
# Row 27
$ws.Range("H27").Value = 3082.8333
$ws.Range("I27").Value = 1820
$ws.Range("J27").Value = 3984.8572
$ws.Range("K27").Value = 1820
$ws.Range("L27").Value = 3984.8572
$ws.Range("M27").Value = -1713
$ws.Range("N27").Value = -4198.8572
# Row 93
$ws.Range("H93").Value = 1864.4
$ws.Range("I93").Value = 1420.4615
$ws.Range("K93").Value = 1420.4615
$ws.Range("M93").Value = -172.4614999999999
# Row 126
$ws.Range("H126").Value = 4574.694
$ws.Range("J126").Value = 4840.143
$ws.Range("L126").Value = 14520.429
$ws.Range("N126").Value = -19460.429
# Row 132
$ws.Range("H132").Value = 222226240
$ws.Range("I132").Value = 3699.75
$ws.Range("J132").Value = 400004300
$ws.Range("K132").Value = 11099.25
$ws.Range("L132").Value = 1200012900
$ws.Range("M132").Value = -8569.25
$ws.Range("N132").Value = -1200017960

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 33
$ws.Range("H33").Value = 5305
$ws.Range("I33").Value = 406.33334
$ws.Range("J33").Value = 20001
$ws.Range("K33").Value = 406.33334
$ws.Range("L33").Value = 20001
$ws.Range("M33").Value = -156.33334
$ws.Range("N33").Value = -20501
# Row 36
$ws.Range("H36").Value = 5305
$ws.Range("I36").Value = 406.33334
$ws.Range("J36").Value = 20001
$ws.Range("K36").Value = 406.33334
$ws.Range("L36").Value = 20001
$ws.Range("M36").Value = -156.33334
$ws.Range("N36").Value = -20501
# Row 38
$ws.Range("H38").Value = 11679501
$ws.Range("I38").Value = 11679501
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 11679501
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -11679028
$ws.Range("N38").Value = ""
# Row 39
$ws.Range("H39").Value = 24998.143
$ws.Range("J39").Value = 24998.143
$ws.Range("L39").Value = 24998.143
$ws.Range("N39").Value = -25824.143
# Row 42
$ws.Range("H42").Value = 10000
$ws.Range("I42").Value = 10000
$ws.Range("K42").Value = 10000
$ws.Range("M42").Value = -9622
# Row 43
$ws.Range("H43").Value = 30000
$ws.Range("I43").Value = 30000
$ws.Range("K43").Value = 30000
$ws.Range("M43").Value = -29851
# Row 44
$ws.Range("H44").Value = 37495
$ws.Range("J44").Value = 37495
$ws.Range("L44").Value = 37495
$ws.Range("N44").Value = -38603
# Row 47
$ws.Range("H47").Value = 35998.668
# Row 52
$ws.Range("H52").Value = 29355.143
$ws.Range("I52").Value = 27398.4
$ws.Range("K52").Value = 27398.4
$ws.Range("M52").Value = -27172.4
# Row 61
$ws.Range("H61").Value = 25747.5
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 25747.5
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 25747.5
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = -26331.5
# Row 62
$ws.Range("H62").Value = 11250.75
$ws.Range("I62").Value = 11250.75
$ws.Range("K62").Value = 11250.75
$ws.Range("M62").Value = -10626.75
# Row 65
$ws.Range("H65").Value = 11250.75
$ws.Range("I65").Value = 11250.75
$ws.Range("K65").Value = 56253.75
$ws.Range("M65").Value = -53133.75
# Row 132
$ws.Range("H132").Value = 8144.636
$ws.Range("I132").Value = 7784.737
$ws.Range("J132").Value = 10424
$ws.Range("K132").Value = 23354.211
$ws.Range("L132").Value = 31272
$ws.Range("M132").Value = -20824.211
$ws.Range("N132").Value = -36332
# Row 136
$ws.Range("H136").Value = 1167.7858
$ws.Range("I136").Value = 1134.6154
$ws.Range("J136").Value = 1599
$ws.Range("K136").Value = 3403.8462
$ws.Range("L136").Value = 4797
$ws.Range("M136").Value = -853.8462
$ws.Range("N136").Value = -9897
# Row 138
$ws.Range("H138").Value = 88750
$ws.Range("J138").Value = 88750
$ws.Range("L138").Value = 88750
$ws.Range("N138").Value = -99030
